$d = $word.ActiveDocument

function DoReplace($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $ok = $find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "FAILED to find/replace: $findText"
    }
}

# --- Overview paragraph: typo fix / rewording ---
DoReplace "in the correct order in order to complete their task effectively." "in the correct order to complete their task corretly."

# --- Key Concepts paragraph: arrays -> dictionaries ---
DoReplace "basic understanding of arrays, consistency" "basic understanding of dictionaries, consistency"

# --- "How you are learning" bullet 1: append read-modify-write routine mention ---
DoReplace "Put a series of transactions in order that makes them atomic." "Put a series of transactions in order that makes them atomic, using read-modify-write routine"

# --- Preparation heading: drop trailing space ---
DoReplace "Preparation " "Preparation"

# --- Bank-transaction bullet: "wrong value"/"green blocks" -> typoed text ---
DoReplace "could hold the wrong value. Show the class the set of green blocks" "could hold an unexcpeted value. Show the class the set of green bloc ks"

# --- "Outline of Activity" heading: drop trailing space before paragraph mark ---
$pOutline = $d.Paragraphs.Item(29)
$endPos = $pOutline.Range.End - 1
$trailing = $d.Range($endPos - 1, $endPos)
if ($trailing.Text -eq " ") {
    $trailing.Text = ""
}

# --- "Ask the students" bullet: single trailing space -> double trailing space ---
$pAsk = $d.Paragraphs.Item(35)
$askEnd = $pAsk.Range.End - 1
$insertPoint = $d.Range($askEnd, $askEnd)
$insertPoint.InsertAfter(" ")

# --- "Learning Intentions" bullets: swap the order of the two items ---
$pFirst = $d.Paragraphs.Item(13)
$pSecond = $d.Paragraphs.Item(14)
$firstText = $pFirst.Range.Text
$secondText = $pSecond.Range.Text
$firstText = $firstText.Substring(0, $firstText.Length - 1)
$secondText = $secondText.Substring(0, $secondText.Length - 1)

$rngFirst = $pFirst.Range
$rngFirst.MoveEnd(1, -1)
$rngFirst.Text = $secondText

$rngSecond = $pSecond.Range
$rngSecond.MoveEnd(1, -1)
$rngSecond.Text = $firstText

Write-Output "Edit complete"
